$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (refresh run)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.444.77'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.697.78'
$ws.Range("E3").Value = '  +0.94%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.94'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5478'
$ws.Range("E6").Value = '  +3.74%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2736'
$ws.Range("E8").Value = '  +1.27%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06447'
$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.01'
$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07704'
$ws.Range("E11").Value = '  +2.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.691.38'
$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.552'
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5844'
$ws.Range("E14").Value = '  +0.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008411'
$ws.Range("E15").Value = '  -1.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.70'
$ws.Range("E16").Value = '  +2.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.475.86'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.949'
$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("E19").Value = '  +0.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.99'
$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.10'
$ws.Range("E21").Value = '  +0.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.259'
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("E23").Value = '  +0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.15'
$ws.Range("E24").Value = '  +3.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1320'
$ws.Range("E25").Value = '  +6.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.894'
$ws.Range("E26").Value = '  +1.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.80'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("E28").Value = '  -5.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.378'
$ws.Range("E29").Value = '  +1.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.332'
$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.613'
$ws.Range("E31").Value = '  +1.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.599'
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.694'
$ws.Range("E33").Value = '  +1.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.041'
$ws.Range("E34").Value = '  +1.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6174'
$ws.Range("E35").Value = '  -0.45%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.411'
$ws.Range("E36").Value = '  +0.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.771'
$ws.Range("E37").Value = '  +2.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01644'
$ws.Range("E38").Value = '  +1.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.118.40'
$ws.Range("E39").Value = '  +0.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.128'
$ws.Range("E40").Value = '  -4.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8788'
$ws.Range("E41").Value = '  +0.20%  '

$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.20'
$ws.Range("E43").Value = '  +0.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.851.13'
$ws.Range("E44").Value = '  +1.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.59'
$ws.Range("E46").Value = '  +1.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.239'
$ws.Range("E47").Value = '  +0.93%  '

$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05285'
$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.140'
$ws.Range("E50").Value = '  +1.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4301'
$ws.Range("E51").Value = '  -0.07%  '
